$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 9 de Agosto de 2020 a las 10:54"

# Updated COVID-19 country statistics
# Row 6
$ws.Range("B6").Value = 2156756
$ws.Range("C6").Value = 4736
$ws.Range("D6").Value = 1481825
$ws.Range("E6").Value = 631433
$ws.Range("G6").Value = 45
$ws.Range("H6").Value = 43498
# Row 25
$ws.Range("B25").Value = 129913
$ws.Range("C25").Value = 3028
$ws.Range("D25").Value = 67673
$ws.Range("E25").Value = 59970
$ws.Range("G25").Value = 61
$ws.Range("H25").Value = 2270
# Row 34
$ws.Range("B34").Value = 82515
$ws.Range("C34").Value = 191
$ws.Range("D34").Value = 57483
$ws.Range("E34").Value = 24435
$ws.Range("G34").Value = 4
$ws.Range("H34").Value = 597
# Row 47
$ws.Range("B47").Value = 55104
$ws.Range("C47").Value = 175
$ws.Range("E47").Value = 6494
# Row 49
$ws.Range("B49").Value = 51791
$ws.Range("C49").Value = 624
$ws.Range("D49").Value = 36691
$ws.Range("E49").Value = 13293
$ws.Range("G49").Value = 7
$ws.Range("H49").Value = 1807
# Row 53
$ws.Range("E53").Value = 2918
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 162
# Row 55
$ws.Range("B55").Value = 40410
$ws.Range("C55").Value = 225
$ws.Range("D55").Value = 32520
$ws.Range("E55").Value = 7099
$ws.Range("G55").Value = 6
$ws.Range("H55").Value = 791
# Row 71
$ws.Range("B71").Value = 22033
$ws.Range("C71").Value = 114
$ws.Range("D71").Value = 19923
$ws.Range("E71").Value = 1389
# Row 111
$ws.Range("B111").Value = 4080
$ws.Range("C111").Value = 72
$ws.Range("D111").Value = 2847
$ws.Range("E111").Value = 1182
# Row 120
$ws.Range("D120").Value = 2579
$ws.Range("E120").Value = 251
# Row 174
$ws.Range("B174").Value = 303
$ws.Range("C174").Value = 8
$ws.Range("D174").Value = 215
$ws.Range("E174").Value = 88
# Row 211
$ws.Range("A211").Value = "San Bartolome"
$ws.Range("C211").Value = 4
$ws.Range("D211").Value = 6
$ws.Range("E211").Value = 7
# Row 212
$ws.Range("A212").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("D212").Value = 7
$ws.Range("E212").Value = 6
$ws.Range("H212").Value = 0
# Row 214
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("H214").Value = 1
# Row 215
$ws.Range("A215").Value = "Santa Sede"
$ws.Range("B215").Value = 12
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("H215").Value = 0
# Row 216
$ws.Range("A216").Value = "Sahara Occidental"
$ws.Range("B216").Value = 10
$ws.Range("D216").Value = 8
$ws.Range("E216").Value = 1
$ws.Range("H216").Value = 1
